# Edits Product Descriptions_NS.docx per commit "made some small edits (mostly punctuation)"
# Strategy: scope Find.Execute to each specific paragraph's Range so that
# repeated phrases (e.g. "No off-putting odor.") are only touched in the
# paragraphs the diff actually changes.

$d = $word.ActiveDocument

function Replace-InPara($index, $find, $replace) {
    $r = $d.Paragraphs($index).Range
    $ok = $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output ("MISS para=$index find=" + $find)
    }
}

# ---------- PLA ----------
Replace-InPara 1 "1,75mm PLA Filament" "1.75mm PLA Filament"
Replace-InPara 2 "The popular and easy to use 3D printer filament, polylactic acid (PLA), available in a wide range of colours." "The popular and easy-to-use 3D printer filament, polylactic acid (PLA), is available in a wide range of colours."
Replace-InPara 10 "Non-mechanical prints like toys and figurines." "Non-mechanical prints like toys and figurines"
Replace-InPara 16 "Easy to use." "Easy to use"
Replace-InPara 17 "No off-putting odor." "No off-putting odor"
Replace-InPara 18 "More environmentally friendly (as compared to other 3D printer filaments)." "More environmentally friendly (as compared to other 3D printer filaments)"
Replace-InPara 20 "Brittle, avoid using for projects that will be bent, twisted, or dropped." "Brittle: avoid using for projects that will be bent, twisted, or dropped"
Replace-InPara 21 "Deforms above temperatures of 60" "Deforms above temperatures of 60"
Replace-InPara 21 "C." "C"

# ---------- PLA (Glow-in-the-dark) ----------
Replace-InPara 23 "1,75mm PLA Filament" "1.75mm PLA Filament"
Replace-InPara 24 "easy to use" "easy-to-use"
Replace-InPara 32 "Halloween projects, wearable prints like jewellery, toys, figurines." "Halloween projects, wearable prints like jewellery, toys, figurines"
Replace-InPara 39 "Easy to use." "Easy to use"
Replace-InPara 40 "No off-putting odor." "No off-putting odor"
Replace-InPara 41 "More environmentally friendly (as compared to other 3D printer filaments)." "More environmentally friendly (as compared to other 3D printer filaments)"
Replace-InPara 43 "Brittle, avoid using for projects that will be bent, twisted, or dropped." "Brittle: avoid using for projects that will be bent, twisted, or dropped"
Replace-InPara 44 "C." "C"

# ---------- ABS ----------
Replace-InPara 46 "1,75mm ABS Filament" "1.75mm ABS Filament"
Replace-InPara 47 "A durable and great material for general purpose projects, acrylonitrile butadiene styrene (ABS) filament is available in a wide range of colours." "A durable and great material for general-purpose projects. Acrylonitrile butadiene styrene (ABS) filament is available in a wide range of colours."
Replace-InPara 55 "Frequently handled projects that may be dropped or heated, like phone cases and electrical enclosures." "Frequently handled projects that may be dropped or heated, like phone cases and electrical enclosures"
Replace-InPara 61 "Superior quality to PLA filament." "Superior quality to PLA filament"
Replace-InPara 62 "No off-putting odor." "No off-putting odor"
Replace-InPara 63 "Strong, durable, and temperature resistant." "Strong, durable, and temperature resistant"
Replace-InPara 65 "Difficult to print." "Difficult to print"
Replace-InPara 66 "Harsh fumes." "Harsh fumes"
Replace-InPara 67 "Prone to warping without the use of a heated bed." "Prone to warping without the use of a heated bed"

# ---------- PETG ----------
Replace-InPara 69 "1,75mm PETG Filament" "1.75mm PETG Filament"
Replace-InPara 70 "A variant of one of the most used plastics in the world, polyethylene terephthalate (PET) filament, available in a wide range of colours." "Polyethylene terephthalate (PET) filament is a variant of one of the most used plastics in the world. It is available in a wide range of colours."
Replace-InPara 78 "Functional objects that may experience physical stress, like mechanical and protective parts." "Functional objects that may experience physical stress, like mechanical and protective parts"
Replace-InPara 82 "Use a low print speed for a higher quality result." "Use a low print speed for a higher quality result"
Replace-InPara 85 "A happy medium between PLA and ABS filaments." "A happy medium between PLA and ABS filaments"
Replace-InPara 86 "Clearer, less brittle, flexible, durable, and temperature resistant." "Clearer, less brittle, flexible, durable, and temperature resistant"
Replace-InPara 87 "Great for layer adhesion." "Great for layer adhesion"
Replace-InPara 89 "Sticky when printed." "Sticky when printed"
Replace-InPara 90 "Scratches more easily." "Scratches easily"
Replace-InPara 91 "Susceptible to moisture." "Susceptible to moisture"

# ---------- TPE ----------
Replace-InPara 93 "1,75mm TPE Flexible Filament" "1.75mm TPE Flexible Filament"
Replace-InPara 102 ", household appliances, medical supplies" ", household appliances, and medical supplies"
Replace-InPara 105 "Not required." "Not required"
Replace-InPara 106 "Tight filament path and slow print speed recommended." "Tight filament path and slow print speed are recommended"
Replace-InPara 109 "Withstands physical stressors that ABS and PLA filaments can" "Withstands physical stressors that ABS and PLA filaments can"
Replace-InPara 109 "tolerate." "tolerate"
Replace-InPara 111 "Can be difficult to extrude." "Can be difficult to extrude"

# ---------- Nylon (PA) ----------
Replace-InPara 113 "1,75mm Nylon (PA) Filament" "1.75mm Nylon (PA) Filament"
Replace-InPara 114 "A popular synthetic polymer, nylon or polyamide (PA) is a go-to filament material for 3D printing. Available in a wide range of colours." "A popular synthetic polymer. Nylon or polyamide (PA) is a go-to filament material for 3D printing and is available in a wide range of colours."
Replace-InPara 122 "mechanical parts like hinges or gears." "mechanical parts like hinges or gears"
Replace-InPara 126 "Use a high nozzle and heated printer bed for best results." "Use a high nozzle and heated printer bed for best results"
Replace-InPara 129 "Can be dyed before or after the printing process." "Can be dyed before or after the printing process"
Replace-InPara 130 "Strong, flexible, durable." "Strong, flexible, durable"
Replace-InPara 132 "Must be stored in a cool, dry place to avoid absorbing moisture." "Must be stored in a cool, dry place to avoid absorbing moisture"

# ---------- Polycarbonate (PC) ----------
Replace-InPara 134 "1,75mm Polycarbonate (PC) Filament" "1.75mm Polycarbonate (PC) Filament"
Replace-InPara 135 " Available in a wide range of colours." " It is available in a wide range of colours."
Replace-InPara 149 "high temperatures and physical stress." "high temperatures and physical stress"
Replace-InPara 151 "Must be stored in a cool, dry place to avoid absorbing moisture." "Must be stored in a cool, dry place to avoid absorbing moisture"
Replace-InPara 152 "Requires a very high print temperature." "Requires a very high print temperature"

# ---------- 3D Pen ----------
Replace-InPara 155 "Create 3D projects without software or files, with a 3D Pen! With a 3D Pen you can draw" "Create 3D projects without software or files! With a 3D Pen you can draw"
Replace-InPara 158 "Pen" "pen"
Replace-InPara 160 ": Use ABS filament material with your 3D Pen, best for beginners and drawing in mid-air." ": Best for beginners and drawing in mid-air"
Replace-InPara 161 ": Use PLA filament material with your 3D Pen, best for drawing directly onto flat surfaces." ": Best for drawing directly onto flat surfaces"
Replace-InPara 162 "Kids, beginners, artists, and educators." "Kids, beginners, artists, and educators"
